$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 17.87134533333333
$ws.Range("H2").Value = 53.614036
$ws.Range("I2").Value = 0.1500697615111392
$ws.Range("J2").Value = 0.1500697615111392
$ws.Range("M2").Value = 35.04689966666667
$ws.Range("N2").Value = 105.140699
$ws.Range("O2").Value = 0.3824629895491901
$ws.Range("P2").Value = 0.3824629895491901
$ws.Range("Q2").Value = 626.335246805685
$ws.Range("R2").Value = 5637.017221251164
$ws.Range("S2").Value = 0.05739612962848428
$ws.Range("T2").Value = 0.05739612962848427
$ws.Range("G3").Value = 17.87134533333333
$ws.Range("H3").Value = 53.614036
$ws.Range("I3").Value = 0.1500697615111392
$ws.Range("J3").Value = 0.1500697615111392
$ws.Range("O3").Value = 0.3264402385872224
$ws.Range("P3").Value = 0.3264402385872223
$ws.Range("Q3").Value = 534.5903603478947
$ws.Range("R3").Value = 4811.313243131051
$ws.Range("S3").Value = 0.04898880875242383
$ws.Range("T3").Value = 0.04898880875242383
$ws.Range("G4").Value = 17.87134533333333
$ws.Range("H4").Value = 53.614036
$ws.Range("I4").Value = 0.1500697615111392
$ws.Range("J4").Value = 0.1500697615111392
$ws.Range("M4").Value = 8.911727666666666
$ws.Range("N4").Value = 26.735183
$ws.Range("O4").Value = 0.09725271102035077
$ws.Range("P4").Value = 0.09725271102035075
$ws.Range("Q4").Value = 159.2645626476209
$ws.Range("R4").Value = 1433.381063828588
$ws.Range("S4").Value = 0.01459469114913578
$ws.Range("T4").Value = 0.01459469114913577
$ws.Range("G5").Value = 17.87134533333333
$ws.Range("H5").Value = 53.614036
$ws.Range("I5").Value = 0.1500697615111392
$ws.Range("J5").Value = 0.1500697615111392
$ws.Range("M5").Value = 17.76285166666667
$ws.Range("N5").Value = 53.288555
$ws.Range("O5").Value = 0.1938440608432367
$ws.Range("P5").Value = 0.1938440608432367
$ws.Range("Q5").Value = 317.4460562397755
$ws.Range("R5").Value = 2857.01450615798
$ws.Range("S5").Value = 0.02909013198109529
$ws.Range("T5").Value = 0.02909013198109529
$ws.Range("I6").Value = 0.2793179663930228
$ws.Range("J6").Value = 0.2793179663930228
$ws.Range("M6").Value = 35.04689966666667
$ws.Range("N6").Value = 105.140699
$ws.Range("O6").Value = 0.3824629895491901
$ws.Range("P6").Value = 0.3824629895491901
$ws.Range("Q6").Value = 1165.769077370395
$ws.Range("R6").Value = 10491.92169633355
$ws.Range("S6").Value = 0.1068287844614757
$ws.Range("T6").Value = 0.1068287844614757
$ws.Range("I7").Value = 0.2793179663930228
$ws.Range("J7").Value = 0.2793179663930228
$ws.Range("O7").Value = 0.3264402385872224
$ws.Range("P7").Value = 0.3264402385872223
$ws.Range("S7").Value = 0.09118062359103611
$ws.Range("T7").Value = 0.0911806235910361
$ws.Range("I8").Value = 0.2793179663930228
$ws.Range("J8").Value = 0.2793179663930228
$ws.Range("M8").Value = 8.911727666666666
$ws.Range("N8").Value = 26.735183
$ws.Range("O8").Value = 0.09725271102035077
$ws.Range("P8").Value = 0.09725271102035075
$ws.Range("Q8").Value = 296.4318281661667
$ws.Range("R8").Value = 2667.886453495501
$ws.Range("S8").Value = 0.02716442946841269
$ws.Range("T8").Value = 0.02716442946841268
$ws.Range("I9").Value = 0.2793179663930228
$ws.Range("J9").Value = 0.2793179663930228
$ws.Range("M9").Value = 17.76285166666667
$ws.Range("N9").Value = 53.288555
$ws.Range("O9").Value = 0.1938440608432367
$ws.Range("P9").Value = 0.1938440608432367
$ws.Range("Q9").Value = 590.847789558176
$ws.Range("R9").Value = 5317.630106023585
$ws.Range("S9").Value = 0.05414412887209826
$ws.Range("T9").Value = 0.05414412887209825
$ws.Range("G10").Value = 12.60542466666667
$ws.Range("H10").Value = 37.816274
$ws.Range("I10").Value = 0.1058506250195358
$ws.Range("J10").Value = 0.1058506250195358
$ws.Range("M10").Value = 35.04689966666667
$ws.Range("N10").Value = 105.140699
$ws.Range("O10").Value = 0.3824629895491901
$ws.Range("P10").Value = 0.3824629895491901
$ws.Range("Q10").Value = 441.7810535483918
$ws.Range("R10").Value = 3976.029481935526
$ws.Range("S10").Value = 0.04048394649062196
$ws.Range("T10").Value = 0.04048394649062196
$ws.Range("G11").Value = 12.60542466666667
$ws.Range("H11").Value = 37.816274
$ws.Range("I11").Value = 0.1058506250195358
$ws.Range("J11").Value = 0.1058506250195358
$ws.Range("O11").Value = 0.3264402385872224
$ws.Range("P11").Value = 0.3264402385872223
$ws.Range("Q11").Value = 377.0694589132353
$ws.Range("R11").Value = 3393.625130219118
$ws.Range("S11").Value = 0.03455390328598387
$ws.Range("T11").Value = 0.03455390328598387
$ws.Range("G12").Value = 12.60542466666667
$ws.Range("H12").Value = 37.816274
$ws.Range("I12").Value = 0.1058506250195358
$ws.Range("J12").Value = 0.1058506250195358
$ws.Range("M12").Value = 8.911727666666666
$ws.Range("N12").Value = 26.735183
$ws.Range("O12").Value = 0.09725271102035077
$ws.Range("P12").Value = 0.09725271102035075
$ws.Range("Q12").Value = 112.3361117520158
$ws.Range("R12").Value = 1011.025005768142
$ws.Range("S12").Value = 0.01029426024634843
$ws.Range("T12").Value = 0.01029426024634843
$ws.Range("G13").Value = 12.60542466666667
$ws.Range("H13").Value = 37.816274
$ws.Range("I13").Value = 0.1058506250195358
$ws.Range("J13").Value = 0.1058506250195358
$ws.Range("M13").Value = 17.76285166666667
$ws.Range("N13").Value = 53.288555
$ws.Range("O13").Value = 0.1938440608432367
$ws.Range("P13").Value = 0.1938440608432367
$ws.Range("Q13").Value = 223.9082885493411
$ws.Range("R13").Value = 2015.17459694407
$ws.Range("S13").Value = 0.02051851499658153
$ws.Range("T13").Value = 0.02051851499658153
$ws.Range("G14").Value = 55.34703199999999
$ws.Range("H14").Value = 166.041096
$ws.Range("I14").Value = 0.4647616470763022
$ws.Range("J14").Value = 0.4647616470763023
$ws.Range("M14").Value = 35.04689966666667
$ws.Range("N14").Value = 105.140699
$ws.Range("O14").Value = 0.3824629895491901
$ws.Range("P14").Value = 0.3824629895491901
$ws.Range("Q14").Value = 1939.741877351789
$ws.Range("R14").Value = 17457.6768961661
$ws.Range("S14").Value = 0.1777541289686081
$ws.Range("T14").Value = 0.1777541289686082
$ws.Range("G15").Value = 55.34703199999999
$ws.Range("H15").Value = 166.041096
$ws.Range("I15").Value = 0.4647616470763022
$ws.Range("J15").Value = 0.4647616470763023
$ws.Range("O15").Value = 0.3264402385872224
$ws.Range("P15").Value = 0.3264402385872223
$ws.Range("Q15").Value = 1655.610656567608
$ws.Range("R15").Value = 14900.49590910847
$ws.Range("S15").Value = 0.1517169029577785
$ws.Range("T15").Value = 0.1517169029577785
$ws.Range("G16").Value = 55.34703199999999
$ws.Range("H16").Value = 166.041096
$ws.Range("I16").Value = 0.4647616470763022
$ws.Range("J16").Value = 0.4647616470763023
$ws.Range("M16").Value = 8.911727666666666
$ws.Range("N16").Value = 26.735183
$ws.Range("O16").Value = 0.09725271102035077
$ws.Range("P16").Value = 0.09725271102035075
$ws.Range("Q16").Value = 493.2376763422852
$ws.Range("R16").Value = 4439.139087080567
$ws.Range("S16").Value = 0.04519933015645387
$ws.Range("T16").Value = 0.04519933015645387
$ws.Range("G17").Value = 55.34703199999999
$ws.Range("H17").Value = 166.041096
$ws.Range("I17").Value = 0.4647616470763022
$ws.Range("J17").Value = 0.4647616470763023
$ws.Range("M17").Value = 17.76285166666667
$ws.Range("N17").Value = 53.288555
$ws.Range("O17").Value = 0.1938440608432367
$ws.Range("P17").Value = 0.1938440608432367
$ws.Range("Q17").Value = 983.1211196062532
$ws.Range("R17").Value = 8848.090076456279
$ws.Range("S17").Value = 0.09009128499346164
$ws.Range("T17").Value = 0.09009128499346165
